# Shift the yearly Coliflor price-series rows down by one and insert a new
# latest reading in row 58 (2021-11-?? -> serial 44519), per the author's
# weekly update. Row 172 is newly created, taking on the data that used to
# live in row 171.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58
$ws.Range("D58").Value = 44519
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 600
$ws.Range("L58").Value = 600
$ws.Range("M58").Value = 600
$ws.Range("P58").Value = 600

# Row 59
$ws.Range("D59").Value = 44392
$ws.Range("J59").Value = 3000
$ws.Range("K59").Value = 700
$ws.Range("L59").Value = 700
$ws.Range("M59").Value = 700
$ws.Range("P59").Value = 700

# Row 60
$ws.Range("D60").Value = 44355
$ws.Range("K60").Value = 450
$ws.Range("L60").Value = 450
$ws.Range("M60").Value = 450
$ws.Range("P60").Value = 450

# Row 61
$ws.Range("D61").Value = 44489
$ws.Range("J61").Value = 4000

# Row 62
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 600
$ws.Range("M62").Value = 600
$ws.Range("P62").Value = 600

# Row 63
$ws.Range("D63").Value = 44434
$ws.Range("I63").Value = "Segunda"
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 500
$ws.Range("L63").Value = 500
$ws.Range("M63").Value = 500
$ws.Range("P63").Value = 500

# Row 64
$ws.Range("D64").Value = 44497
$ws.Range("J64").Value = 5000

# Row 65
$ws.Range("D65").Value = 44449
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 600
$ws.Range("L65").Value = 600
$ws.Range("M65").Value = 600
$ws.Range("P65").Value = 600

# Row 66
$ws.Range("D66").Value = 44358
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 5000

# Row 67
$ws.Range("D67").Value = 44399
$ws.Range("I67").Value = "Segunda"
$ws.Range("K67").Value = 500
$ws.Range("L67").Value = 500
$ws.Range("M67").Value = 500
$ws.Range("P67").Value = 500

# Row 68
$ws.Range("D68").Value = 44298
$ws.Range("I68").Value = "Primera"
$ws.Range("K68").Value = 700
$ws.Range("L68").Value = 700
$ws.Range("M68").Value = 700
$ws.Range("P68").Value = 700

# Row 69
$ws.Range("D69").Value = 44405
$ws.Range("I69").Value = "Segunda"
$ws.Range("K69").Value = 500
$ws.Range("L69").Value = 500
$ws.Range("M69").Value = 500
$ws.Range("P69").Value = 500

# Row 70
$ws.Range("D70").Value = 44273
$ws.Range("K70").Value = 800
$ws.Range("L70").Value = 800
$ws.Range("M70").Value = 800
$ws.Range("O70").Value = "Región del Maule"
$ws.Range("P70").Value = 800

# Row 71
$ws.Range("J71").Value = 3000
$ws.Range("L71").Value = 600
$ws.Range("M71").Value = 600
$ws.Range("O71").Value = "Región Metropolitana"
$ws.Range("P71").Value = 600

# Row 72
$ws.Range("I72").Value = "Primera"
$ws.Range("K72").Value = 600
$ws.Range("L72").Value = 650
$ws.Range("M72").Value = 612
$ws.Range("P72").Value = 612

# Row 73
$ws.Range("D73").Value = 44435
$ws.Range("I73").Value = "Segunda"
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 500
$ws.Range("L73").Value = 500
$ws.Range("M73").Value = 500
$ws.Range("P73").Value = 500

# Row 74
$ws.Range("D74").Value = 44328
$ws.Range("J74").Value = 300
$ws.Range("K74").Value = 700
$ws.Range("L74").Value = 700
$ws.Range("M74").Value = 700
$ws.Range("P74").Value = 700

# Row 75
$ws.Range("D75").Value = 44277
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 800
$ws.Range("L75").Value = 800
$ws.Range("M75").Value = 800
$ws.Range("P75").Value = 800

# Row 76
$ws.Range("D76").Value = 44516
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 500
$ws.Range("L76").Value = 500
$ws.Range("M76").Value = 500
$ws.Range("P76").Value = 500

# Row 77
$ws.Range("D77").Value = 44168

# Row 78
$ws.Range("D78").Value = 44475
$ws.Range("I78").Value = "Primera"
$ws.Range("K78").Value = 600
$ws.Range("L78").Value = 600
$ws.Range("M78").Value = 600
$ws.Range("P78").Value = 600

# Row 79
$ws.Range("D79").Value = 44419
$ws.Range("I79").Value = "Segunda"
$ws.Range("J79").Value = 3000

# Row 80
$ws.Range("D80").Value = 44162
$ws.Range("J80").Value = 5000

# Row 81
$ws.Range("D81").Value = 44357
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 500
$ws.Range("L81").Value = 500
$ws.Range("M81").Value = 500
$ws.Range("P81").Value = 500

# Row 82
$ws.Range("D82").Value = 44333
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 600
$ws.Range("L82").Value = 600
$ws.Range("M82").Value = 600
$ws.Range("P82").Value = 600

# Row 83
$ws.Range("D83").Value = 44320
$ws.Range("K83").Value = 700
$ws.Range("L83").Value = 700
$ws.Range("M83").Value = 700
$ws.Range("P83").Value = 700

# Row 84
$ws.Range("D84").Value = 44467
$ws.Range("K84").Value = 600
$ws.Range("L84").Value = 600
$ws.Range("M84").Value = 600
$ws.Range("P84").Value = 600

# Row 85
$ws.Range("D85").Value = 44264
$ws.Range("K85").Value = 800
$ws.Range("L85").Value = 800
$ws.Range("M85").Value = 800
$ws.Range("P85").Value = 800

# Row 86
$ws.Range("D86").Value = 44214
$ws.Range("K86").Value = 700
$ws.Range("L86").Value = 700
$ws.Range("M86").Value = 700
$ws.Range("P86").Value = 700

# Row 87
$ws.Range("D87").Value = 44167
$ws.Range("J87").Value = 3000
$ws.Range("K87").Value = 600
$ws.Range("L87").Value = 600
$ws.Range("M87").Value = 600
$ws.Range("P87").Value = 600

# Row 88
$ws.Range("D88").Value = 44291
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 700
$ws.Range("L88").Value = 700
$ws.Range("M88").Value = 700
$ws.Range("P88").Value = 700

# Row 89
$ws.Range("D89").Value = 44293
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 600
$ws.Range("L89").Value = 600
$ws.Range("M89").Value = 600
$ws.Range("P89").Value = 600

# Row 90
$ws.Range("D90").Value = 44496
$ws.Range("J90").Value = 5000

# Row 91
$ws.Range("D91").Value = 44326
$ws.Range("K91").Value = 500
$ws.Range("L91").Value = 500
$ws.Range("M91").Value = 500
$ws.Range("P91").Value = 500

# Row 92
$ws.Range("D92").Value = 44302
$ws.Range("K92").Value = 600
$ws.Range("L92").Value = 600
$ws.Range("M92").Value = 600
$ws.Range("P92").Value = 600

# Row 93
$ws.Range("D93").Value = 44292
$ws.Range("K93").Value = 700
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = 700
$ws.Range("P93").Value = 700

# Row 94
$ws.Range("D94").Value = 44308
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 600
$ws.Range("L94").Value = 600
$ws.Range("M94").Value = 600
$ws.Range("P94").Value = 600

# Row 95
$ws.Range("D95").Value = 44498
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 5000
$ws.Range("K95").Value = 550
$ws.Range("L95").Value = 550
$ws.Range("M95").Value = 550
$ws.Range("P95").Value = 550

# Row 96
$ws.Range("D96").Value = 44420
$ws.Range("I96").Value = "Segunda"
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 500
$ws.Range("L96").Value = 500
$ws.Range("M96").Value = 500
$ws.Range("O96").Value = "Región del Maule"
$ws.Range("P96").Value = 500

# Row 97
$ws.Range("D97").Value = 44396
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 750
$ws.Range("L97").Value = 750
$ws.Range("M97").Value = 750
$ws.Range("O97").Value = "Región Metropolitana"
$ws.Range("P97").Value = 750

# Row 98
$ws.Range("D98").Value = 44321
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 600
$ws.Range("L98").Value = 600
$ws.Range("M98").Value = 600
$ws.Range("P98").Value = 600

# Row 99
$ws.Range("D99").Value = 44349
$ws.Range("J99").Value = 6000
$ws.Range("K99").Value = 500
$ws.Range("L99").Value = 500
$ws.Range("M99").Value = 500
$ws.Range("P99").Value = 500

# Row 100
$ws.Range("D100").Value = 44477
$ws.Range("J100").Value = 3000

# Row 101
$ws.Range("D101").Value = 44487
$ws.Range("J101").Value = 4000

# Row 102
$ws.Range("D102").Value = 44452
$ws.Range("J102").Value = 3000
$ws.Range("K102").Value = 600
$ws.Range("L102").Value = 600
$ws.Range("M102").Value = 600
$ws.Range("P102").Value = 600

# Row 103
$ws.Range("D103").Value = 44505
$ws.Range("J103").Value = 6000
$ws.Range("K103").Value = 500
$ws.Range("L103").Value = 500
$ws.Range("M103").Value = 500
$ws.Range("P103").Value = 500

# Row 104
$ws.Range("D104").Value = 44306
$ws.Range("J104").Value = 4000

# Row 105
$ws.Range("D105").Value = 44509
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 600
$ws.Range("L105").Value = 600
$ws.Range("M105").Value = 600
$ws.Range("P105").Value = 600

# Row 106
$ws.Range("D106").Value = 44189
$ws.Range("K106").Value = 500
$ws.Range("L106").Value = 500
$ws.Range("M106").Value = 500
$ws.Range("P106").Value = 500

# Row 107
$ws.Range("D107").Value = 44278

# Row 108
$ws.Range("D108").Value = 44265
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 800
$ws.Range("L108").Value = 800
$ws.Range("M108").Value = 800
$ws.Range("P108").Value = 800

# Row 109
$ws.Range("D109").Value = 44494
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 600
$ws.Range("L109").Value = 600
$ws.Range("M109").Value = 600
$ws.Range("P109").Value = 600

# Row 110
$ws.Range("D110").Value = 44300
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 700
$ws.Range("L110").Value = 700
$ws.Range("M110").Value = 700
$ws.Range("P110").Value = 700

# Row 111
$ws.Range("D111").Value = 44356
$ws.Range("J111").Value = 5000
$ws.Range("K111").Value = 450
$ws.Range("L111").Value = 450
$ws.Range("M111").Value = 450
$ws.Range("P111").Value = 450

# Row 112
$ws.Range("D112").Value = 44469
$ws.Range("J112").Value = 4000
$ws.Range("K112").Value = 500
$ws.Range("L112").Value = 500
$ws.Range("M112").Value = 500
$ws.Range("P112").Value = 500

# Row 113
$ws.Range("D113").Value = 44453
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 600
$ws.Range("P113").Value = 600

# Row 114
$ws.Range("D114").Value = 44518
$ws.Range("J114").Value = 4000
$ws.Range("K114").Value = 500
$ws.Range("L114").Value = 500
$ws.Range("M114").Value = 500
$ws.Range("P114").Value = 500

# Row 115
$ws.Range("D115").Value = 44446

# Row 116
$ws.Range("D116").Value = 44463

# Row 117
$ws.Range("D117").Value = 44323
$ws.Range("I117").Value = "Primera"
$ws.Range("K117").Value = 600
$ws.Range("L117").Value = 600
$ws.Range("M117").Value = 600
$ws.Range("P117").Value = 600

# Row 118
$ws.Range("D118").Value = 44417
$ws.Range("I118").Value = "Segunda"
$ws.Range("J118").Value = 3000
$ws.Range("K118").Value = 500
$ws.Range("L118").Value = 500
$ws.Range("M118").Value = 500
$ws.Range("P118").Value = 500

# Row 119
$ws.Range("D119").Value = 44445
$ws.Range("K119").Value = 600
$ws.Range("L119").Value = 600
$ws.Range("M119").Value = 600
$ws.Range("P119").Value = 600

# Row 120
$ws.Range("D120").Value = 44342
$ws.Range("I120").Value = "Primera"
$ws.Range("J120").Value = 4000
$ws.Range("K120").Value = 700
$ws.Range("L120").Value = 700
$ws.Range("M120").Value = 700
$ws.Range("P120").Value = 700

# Row 121
$ws.Range("D121").Value = 44406
$ws.Range("I121").Value = "Segunda"
$ws.Range("J121").Value = 5000
$ws.Range("K121").Value = 400
$ws.Range("L121").Value = 400
$ws.Range("M121").Value = 400
$ws.Range("P121").Value = 400

# Row 122
$ws.Range("D122").Value = 44295
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 700
$ws.Range("L122").Value = 700
$ws.Range("M122").Value = 700
$ws.Range("P122").Value = 700

# Row 123
$ws.Range("D123").Value = 44270
$ws.Range("J123").Value = 3000
$ws.Range("K123").Value = 800
$ws.Range("L123").Value = 800
$ws.Range("M123").Value = 800
$ws.Range("P123").Value = 800

# Row 124
$ws.Range("D124").Value = 44363
$ws.Range("J124").Value = 6000
$ws.Range("K124").Value = 400
$ws.Range("L124").Value = 400
$ws.Range("M124").Value = 400
$ws.Range("P124").Value = 400

# Row 125
$ws.Range("D125").Value = 44299
$ws.Range("J125").Value = 4000
$ws.Range("K125").Value = 700
$ws.Range("L125").Value = 700
$ws.Range("M125").Value = 700
$ws.Range("P125").Value = 700

# Row 126
$ws.Range("D126").Value = 44336
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 600
$ws.Range("L126").Value = 600
$ws.Range("M126").Value = 600
$ws.Range("P126").Value = 600

# Row 127
$ws.Range("D127").Value = 44372
$ws.Range("J127").Value = 5000
$ws.Range("K127").Value = 400
$ws.Range("L127").Value = 400
$ws.Range("M127").Value = 400
$ws.Range("O127").Value = "Región del Maule"
$ws.Range("P127").Value = 400

# Row 128
$ws.Range("I128").Value = "Primera"
$ws.Range("J128").Value = 3000
$ws.Range("K128").Value = 700
$ws.Range("L128").Value = 700
$ws.Range("M128").Value = 700
$ws.Range("O128").Value = "Región Metropolitana"
$ws.Range("P128").Value = 700

# Row 129
$ws.Range("D129").Value = 44403
$ws.Range("I129").Value = "Segunda"
$ws.Range("J129").Value = 2000
$ws.Range("K129").Value = 500
$ws.Range("L129").Value = 500
$ws.Range("M129").Value = 500
$ws.Range("P129").Value = 500

# Row 130
$ws.Range("D130").Value = 44169
$ws.Range("J130").Value = 3000
$ws.Range("K130").Value = 600
$ws.Range("L130").Value = 600
$ws.Range("M130").Value = 600
$ws.Range("P130").Value = 600

# Row 131
$ws.Range("D131").Value = 44376
$ws.Range("J131").Value = 5000
$ws.Range("K131").Value = 400
$ws.Range("L131").Value = 400
$ws.Range("M131").Value = 400
$ws.Range("P131").Value = 400

# Row 132
$ws.Range("D132").Value = 44172
$ws.Range("I132").Value = "Primera"
$ws.Range("K132").Value = 600
$ws.Range("L132").Value = 600
$ws.Range("M132").Value = 600
$ws.Range("P132").Value = 600

# Row 133
$ws.Range("D133").Value = 44421
$ws.Range("I133").Value = "Segunda"
$ws.Range("K133").Value = 500
$ws.Range("L133").Value = 500
$ws.Range("M133").Value = 500
$ws.Range("O133").Value = "Región del Maule"
$ws.Range("P133").Value = 500

# Row 134
$ws.Range("D134").Value = 44431
$ws.Range("I134").Value = "Primera"
$ws.Range("K134").Value = 600
$ws.Range("L134").Value = 600
$ws.Range("M134").Value = 600
$ws.Range("O134").Value = "Región Metropolitana"
$ws.Range("P134").Value = 600

# Row 135
$ws.Range("D135").Value = 44426
$ws.Range("I135").Value = "Segunda"
$ws.Range("K135").Value = 500
$ws.Range("L135").Value = 500
$ws.Range("M135").Value = 500
$ws.Range("P135").Value = 500

# Row 136
$ws.Range("D136").Value = 44448
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 600
$ws.Range("L136").Value = 600
$ws.Range("M136").Value = 600
$ws.Range("P136").Value = 600

# Row 137
$ws.Range("D137").Value = 44362
$ws.Range("J137").Value = 5000
$ws.Range("K137").Value = 400
$ws.Range("L137").Value = 400
$ws.Range("M137").Value = 400
$ws.Range("P137").Value = 400

# Row 138
$ws.Range("D138").Value = 44176
$ws.Range("J138").Value = 3000

# Row 139
$ws.Range("D139").Value = 44301
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 600
$ws.Range("L139").Value = 600
$ws.Range("M139").Value = 600
$ws.Range("P139").Value = 600

# Row 140
$ws.Range("D140").Value = 44407
$ws.Range("I140").Value = "Segunda"
$ws.Range("J140").Value = 3000
$ws.Range("K140").Value = 500
$ws.Range("L140").Value = 500
$ws.Range("M140").Value = 500
$ws.Range("P140").Value = 500

# Row 141
$ws.Range("D141").Value = 44284
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 700
$ws.Range("L141").Value = 800
$ws.Range("M141").Value = 740
$ws.Range("P141").Value = 740

# Row 142
$ws.Range("D142").Value = 44441
$ws.Range("K142").Value = 600
$ws.Range("L142").Value = 600
$ws.Range("M142").Value = 600
$ws.Range("P142").Value = 600

# Row 143
$ws.Range("D143").Value = 44279
$ws.Range("K143").Value = 800
$ws.Range("L143").Value = 800
$ws.Range("M143").Value = 800
$ws.Range("P143").Value = 800

# Row 144
$ws.Range("D144").Value = 44341
$ws.Range("J144").Value = 3000
$ws.Range("K144").Value = 700
$ws.Range("L144").Value = 700
$ws.Range("M144").Value = 700
$ws.Range("P144").Value = 700

# Row 145
$ws.Range("D145").Value = 44504
$ws.Range("K145").Value = 600
$ws.Range("L145").Value = 600
$ws.Range("M145").Value = 600
$ws.Range("P145").Value = 600

# Row 146
$ws.Range("D146").Value = 44350
$ws.Range("K146").Value = 500
$ws.Range("L146").Value = 500
$ws.Range("M146").Value = 500
$ws.Range("O146").Value = "Región del Maule"
$ws.Range("P146").Value = 500

# Row 147
$ws.Range("D147").Value = 44312
$ws.Range("J147").Value = 5000
$ws.Range("O147").Value = "Provincia de Chacabuco"

# Row 148
$ws.Range("D148").Value = 44384
$ws.Range("J148").Value = 4000
$ws.Range("K148").Value = 600
$ws.Range("L148").Value = 600
$ws.Range("M148").Value = 600
$ws.Range("P148").Value = 600

# Row 149
$ws.Range("D149").Value = 44329
$ws.Range("K149").Value = 650
$ws.Range("L149").Value = 650
$ws.Range("M149").Value = 650
$ws.Range("P149").Value = 650

# Row 150
$ws.Range("D150").Value = 44491
$ws.Range("K150").Value = 700
$ws.Range("L150").Value = 700
$ws.Range("M150").Value = 700
$ws.Range("P150").Value = 700

# Row 151
$ws.Range("D151").Value = 44272
$ws.Range("J151").Value = 3000
$ws.Range("K151").Value = 800
$ws.Range("L151").Value = 800
$ws.Range("M151").Value = 800
$ws.Range("P151").Value = 800

# Row 152
$ws.Range("D152").Value = 44305
$ws.Range("J152").Value = 4000
$ws.Range("K152").Value = 600
$ws.Range("L152").Value = 600
$ws.Range("M152").Value = 600
$ws.Range("P152").Value = 600

# Row 153
$ws.Range("D153").Value = 44166
$ws.Range("K153").Value = 700
$ws.Range("L153").Value = 700
$ws.Range("M153").Value = 700
$ws.Range("P153").Value = 700

# Row 154
$ws.Range("D154").Value = 44315
$ws.Range("J154").Value = 3000
$ws.Range("K154").Value = 600
$ws.Range("L154").Value = 600
$ws.Range("M154").Value = 600
$ws.Range("P154").Value = 600

# Row 155
$ws.Range("D155").Value = 44348
$ws.Range("J155").Value = 6000
$ws.Range("K155").Value = 450
$ws.Range("L155").Value = 450
$ws.Range("M155").Value = 450
$ws.Range("P155").Value = 450

# Row 156
$ws.Range("D156").Value = 44322
$ws.Range("J156").Value = 3000
$ws.Range("K156").Value = 600
$ws.Range("L156").Value = 600
$ws.Range("M156").Value = 600
$ws.Range("P156").Value = 600

# Row 157
$ws.Range("D157").Value = 44495
$ws.Range("J157").Value = 5000
$ws.Range("K157").Value = 500
$ws.Range("L157").Value = 500
$ws.Range("M157").Value = 500
$ws.Range("P157").Value = 500

# Row 158
$ws.Range("D158").Value = 44327
$ws.Range("J158").Value = 3000
$ws.Range("K158").Value = 700
$ws.Range("L158").Value = 700
$ws.Range("M158").Value = 700
$ws.Range("P158").Value = 700

# Row 159
$ws.Range("D159").Value = 44510
$ws.Range("K159").Value = 600
$ws.Range("L159").Value = 600
$ws.Range("M159").Value = 600
$ws.Range("P159").Value = 600

# Row 160
$ws.Range("D160").Value = 44161
$ws.Range("J160").Value = 4000

# Row 161
$ws.Range("D161").Value = 44517
$ws.Range("J161").Value = 5000
$ws.Range("K161").Value = 500
$ws.Range("L161").Value = 500
$ws.Range("M161").Value = 500
$ws.Range("P161").Value = 500

# Row 162
$ws.Range("D162").Value = 44391
$ws.Range("K162").Value = 700
$ws.Range("L162").Value = 700
$ws.Range("M162").Value = 700
$ws.Range("P162").Value = 700

# Row 163
$ws.Range("D163").Value = 44340
$ws.Range("J163").Value = 3000
$ws.Range("K163").Value = 600
$ws.Range("L163").Value = 600
$ws.Range("M163").Value = 600
$ws.Range("P163").Value = 600

# Row 164
$ws.Range("D164").Value = 44515
$ws.Range("K164").Value = 500
$ws.Range("L164").Value = 500
$ws.Range("M164").Value = 500
$ws.Range("P164").Value = 500

# Row 165
$ws.Range("D165").Value = 44330
$ws.Range("I165").Value = "Primera"
$ws.Range("J165").Value = 4000
$ws.Range("K165").Value = 650
$ws.Range("L165").Value = 650
$ws.Range("M165").Value = 650
$ws.Range("P165").Value = 650

# Row 166
$ws.Range("D166").Value = 44432
$ws.Range("I166").Value = "Segunda"
$ws.Range("J166").Value = 3000
$ws.Range("K166").Value = 500
$ws.Range("L166").Value = 500
$ws.Range("M166").Value = 500
$ws.Range("P166").Value = 500

# Row 167
$ws.Range("D167").Value = 44181
$ws.Range("J167").Value = 2000

# Row 168
$ws.Range("D168").Value = 44271
$ws.Range("J168").Value = 3000
$ws.Range("K168").Value = 800
$ws.Range("L168").Value = 800
$ws.Range("M168").Value = 800
$ws.Range("P168").Value = 800

# Row 169
$ws.Range("D169").Value = 44307
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 5000

# Row 170
$ws.Range("D170").Value = 44400
$ws.Range("I170").Value = "Segunda"
$ws.Range("J170").Value = 3000
$ws.Range("K170").Value = 500
$ws.Range("L170").Value = 500
$ws.Range("M170").Value = 500
$ws.Range("P170").Value = 500

# Row 171
$ws.Range("D171").Value = 44309
$ws.Range("K171").Value = 600
$ws.Range("L171").Value = 600
$ws.Range("M171").Value = 600
$ws.Range("P171").Value = 600

# Row 172 (new row, content matches the prior row 171 before the shift)
$ws.Range("A172").Value = 5
$ws.Range("B172").Value = "Macroferia Regional de Talca"
$ws.Range("C172").Value = "Maule"
$ws.Range("D172").Value = 44508
$ws.Range("E172").Value = 7
$ws.Range("F172").Value = 100112008
$ws.Range("G172").Value = "Coliflor"
$ws.Range("H172").Value = "Sin especificar"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 5000
$ws.Range("K172").Value = 500
$ws.Range("L172").Value = 500
$ws.Range("M172").Value = 500
$ws.Range("N172").Value = "$/unidad"
$ws.Range("O172").Value = "Región del Maule"
$ws.Range("P172").Value = 500
$ws.Range("Q172").Value = 1
$ws.Range("R172").Value = "Hortaliza"

# Copy the date style (custom date format) from D171 to the new D172 cell
$ws.Range("D171").Copy()
$ws.Range("D172").PasteSpecial(-4122) | Out-Null
$ws.Range("D172").Value = 44508
$excel.CutCopyMode = 0
